# Refresh the cached text of the auto "Date" placeholder (the
# `datetime2`-formatted <a:fld>) that appears on the slide master and on
# every slide layout, from "Thursday, June 5, 2025" to "Tuesday, May 6,
# 2025" - mirroring what PowerPoint does when it recomputes/re-saves an
# auto-updating date field on a later day.
$p = $ppt.ActivePresentation

$oldText = "Thursday, June 5, 2025"
$newText = "Tuesday, May 6, 2025"

function Update-DateShape($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*" -and $shp.HasTextFrame) {
            if ($shp.TextFrame.TextRange.Text -eq $oldText) {
                $shp.TextFrame.TextRange.Text = $newText
            }
        }
    }
}

# Slide master's own Date Placeholder.
Update-DateShape $p.SlideMaster.Shapes

# Every slide layout (CustomLayout) owned by the slide master.
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Update-DateShape $layout.Shapes
}
